$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update model name headers
$ws.Range("C1").Value = "newModel"
$ws.Range("D1").Value = "newModel2"

# Row 3: rename feature LYMC -> LYMF
$ws.Range("A3").Value = "LYMF"

# Row 5: move the value from D5 to C5 (enable prediction through own model)
$ws.Range("D5").ClearContents()
$ws.Range("C5").Value = 1

# Row 6: enable own-model prediction flag
$ws.Range("D6").Value = 1

# Row 8: enable own-model prediction flag
$ws.Range("D8").Value = 1

# Row 9: enable own-model prediction flag
$ws.Range("D9").Value = 1
